$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the existing (pre-edit) formatting of column C onto the new column D,
#    before we touch the font on A:C, so D keeps the original theme font.
$ws.Range("C1:C3").Copy()
$ws.Range("D1:D3").PasteSpecial(-4122)

# 2) Write the new augmented-matrix values (Gaussian elimination matrix).
$ws.Range("A1").Value = 2
$ws.Range("B1").Value = 3
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = -1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2

# 3) The existing matrix columns (A:C) pick up an explicit black font color.
$ws.Range("A1:C3").Font.Color = 0

# 4) Column D matches the width/format of the other columns.
$ws.Columns("D").ColumnWidth = $ws.Columns("A").ColumnWidth
